# PC Value Tracker - Monthly_Report_2025-11.xlsx edit
# Reflects the change from a "Diagnostic" stream issue being reclassified as
# "Day-to-Day" (with its business impact relabeled Compliance -> Efficiency,
# and the PLC issue's Production impact relabeled -> Low), which then
# collapses the "BY STREAM" / "By Stream" breakdowns onto a single stream.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Summary"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Summary")

# Updated "Generated:" timestamp
$ws1.Range("B3").Value = "2026-01-24 21:19"

# The second issue (previously "Diagnostic" stream) is now also
# "Day-to-Day", so the BY STREAM breakdown collapses to one row, and every
# row from "BY SYSTEM" onward shifts up by one.
$ws1.Rows.Item(9).Delete()

# BY STREAM count: both issues are Day-to-Day now
$ws1.Range("B8").Value = 2

# Re-bold the section headers (col A) that are now the first row of each
# block, and un-bold the data rows that used to be bold section headers
# before the shift.
$ws1.Range("A1").Copy()
$ws1.Range("A7").PasteSpecial(-4122)
$ws1.Range("A10").PasteSpecial(-4122)
$ws1.Range("A14").PasteSpecial(-4122)
$ws1.Range("A17").PasteSpecial(-4122)
$ws1.Range("A21").PasteSpecial(-4122)

$ws1.Range("A2").Copy()
$ws1.Range("A8").PasteSpecial(-4122)
$ws1.Range("A12").PasteSpecial(-4122)

# Business impact labels: Production -> Low, Compliance -> Efficiency
$ws1.Range("A18").Value = "Low"
$ws1.Range("A19").Value = "Efficiency"

# ---------------------------------------------------------------------
# Sheet 2: "All Issues"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("All Issues")

$ws2.Range("G2").Value = "Low"
$ws2.Range("D3").Value = "Day-to-Day"
$ws2.Range("G3").Value = "Efficiency"

# ---------------------------------------------------------------------
# Sheet 3: "By Stream"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("By Stream")

# First (Day-to-Day) block: impact relabel
$ws3.Range("F3").Value = "Low"

# Second issue now belongs to the Day-to-Day block too, so it becomes the
# second data row right under the first, and the old "DIAGNOSTIC" block
# (header + its own Date/System/.../Impact row + data row + trailing blank
# row) goes away.
# (Force the date-looking text to stay text instead of being parsed into a
# date serial by first marking it as Text, then restoring the cell's
# formatting to plain/default via a format-only paste from a neighboring
# default-styled cell.)
$ws3.Range("A4").NumberFormat = "@"
$ws3.Range("A4").Value = "2025-11-17"
$ws3.Range("B1").Copy()
$ws3.Range("A4").PasteSpecial(-4122)
$ws3.Range("B4").Value = "Other"
$ws3.Range("C4").Value = "Coordinated restore path for corrupted server pgwgen004002; evaluated backup options"
$ws3.Range("D4").Value = "Moderate"
$ws3.Range("E4").Value = "Pending"
$ws3.Range("F4").Value = "Efficiency"

# Clear the old " === DIAGNOSTIC ===" banner text - row 5 becomes the new
# trailing blank row. Re-paste the (blank) formatting from B5 afterwards so
# the now-empty A5 cell stays present in the sheet (matching its blank
# siblings B5:F5) instead of being dropped entirely.
$ws3.Range("A5").Value = ""
$ws3.Range("B5").Copy()
$ws3.Range("A5").PasteSpecial(-4122)

# Remove the old DIAGNOSTIC header/data/blank rows (previously rows 6-8).
$ws3.Rows("6:8").Delete()

# ---------------------------------------------------------------------
# Sheet 4: "By System"
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("By System")

$ws4.Range("D3").Value = "Day-to-Day"
